$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C38").Value = 87.5
$ws.Range("C39").Value = 87.47
$ws.Range("C40").Value = 88.34999999999999
$ws.Range("C41").Value = 89.31
$ws.Range("C42").Value = 86.63
$ws.Range("C44").Value = 92.22
$ws.Range("C45").Value = 91.48999999999999
$ws.Range("C46").Value = 92.98999999999999
$ws.Range("C48").Value = 91.63
$ws.Range("C49").Value = 95.33
$ws.Range("C50").Value = 95.23
$ws.Range("C51").Value = 96.7
$ws.Range("C52").Value = 96.18000000000001
$ws.Range("C53").Value = 97.54000000000001
$ws.Range("C54").Value = 100.6
$ws.Range("C55").Value = 97.81
$ws.Range("C56").Value = 100.8
$ws.Range("C57").Value = 100.54
$ws.Range("C59").Value = 101.57
$ws.Range("C60").Value = 99.08
$ws.Range("C61").Value = 99.59
$ws.Range("C63").Value = 102.26
$ws.Range("C64").Value = 98.83
$ws.Range("C66").Value = 100.8
$ws.Range("C67").Value = 99.12
$ws.Range("C68").Value = 98.37
$ws.Range("C70").Value = 94.75
$ws.Range("C71").Value = 97.68000000000001
$ws.Range("E71").Value = 95.72
$ws.Range("C72").Value = 100.92
$ws.Range("C74").Value = 102.51
$ws.Range("C75").Value = 102.16
$ws.Range("C76").Value = 101.27
$ws.Range("C77").Value = 102.76
$ws.Range("C78").Value = 100.18
$ws.Range("C79").Value = 102.22
$ws.Range("C80").Value = 101.96
$ws.Range("E80").Value = 97.45999999999999
$ws.Range("E81").Value = 100.31
$ws.Range("C82").Value = 101.83
$ws.Range("C83").Value = 97.92
$ws.Range("E83").Value = 91.89
$ws.Range("C84").Value = 98.94
$ws.Range("C85").Value = 101.07
$ws.Range("C86").Value = 102.15
$ws.Range("B87").Value = 100.74
$ws.Range("C87").Value = 103.03
$ws.Range("E87").Value = 102.56
